# Adds a new "SLA" column (D) with two elapsed-time interval cells and one
# date/time cell, matching the commit "Added interval support, fixed skipping files".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---
$ws.Range("D1").Value = "SLA"

# --- New interval / date values in column D ---
# Apply the [h]:mm:ss (numFmtId 46) style first so it is minted as cellXfs
# index 2, then the m/d/yy h:mm (numFmtId 22) style as cellXfs index 3 -
# matching the style order the workbook ends up with.
$ws.Range("D4").Value = 2.7919328703703705
$ws.Range("D4").NumberFormat = "[h]:mm:ss"

$ws.Range("D2").Value = 1.7017708333333332
$ws.Range("D2").NumberFormat = "m/d/yy h:mm"

$ws.Range("D3").Value = 1.3333333333333333
$ws.Range("D3").NumberFormat = "m/d/yy h:mm"

# --- Column width for the new column ---
$ws.Columns.Item(4).ColumnWidth = 24.6

# --- Selection moves to the newly populated cell D2 ---
[void]$ws.Range("D2").Select()
